# Insert a brand-new "A06 vie saint gregoire" worksheet right after
# "A05 vie sainte dieudonnee", shifting every following sheet's leading
# number up by one (A06->A07, A07->A08, ..., A13->A14, A16->A17, ...,
# A25->A26), and give the new sheet the same 4-column header row used by
# every other sheet (line_n / prev_line / line / next_line), bold +
# centered, with no data rows below it.

$wb = $excel.ActiveWorkbook

# Keep a handle to whichever sheet is active right now so we can restore the
# original active-tab selection once we're done (Add()/rename shouldn't
# change which tab is shown).
$originallyActive = $wb.ActiveSheet

# 1) Rename every sheet from index 6 to the end (1-based), incrementing the
#    leading "A##" number by one. Walk from the last sheet back to index 6
#    so the computed names never collide with a not-yet-renamed sheet.
$count = $wb.Worksheets.Count
for ($i = $count; $i -ge 6; $i--) {
    $ws = $wb.Worksheets.Item($i)
    $oldName = $ws.Name
    $prefixNum = [int]$oldName.Substring(1, 2)
    $rest = $oldName.Substring(3)
    $newNum = $prefixNum + 1
    $newPrefix = "A" + $newNum.ToString("D2")
    $ws.Name = $newPrefix + $rest
}

# 2) Insert the new worksheet right after "A05 vie sainte dieudonnee"
#    (which is still at index 5) so it lands at index 6.
$afterSheet = $wb.Worksheets.Item(5)
$newSheet = $wb.Worksheets.Add($null, $afterSheet)
$newSheet.Name = "A06 vie saint gregoire"

# 3) Populate its header row exactly like the other story sheets.
$newSheet.Range("A1").Value = "line_n"
$newSheet.Range("B1").Value = "prev_line"
$newSheet.Range("C1").Value = "line"
$newSheet.Range("D1").Value = "next_line"
$newSheet.Range("A1:D1").Font.Bold = $true
$newSheet.Range("A1:D1").HorizontalAlignment = -4108

# 4) Restore the original active tab (Worksheets.Add() switches focus to the
#    sheet it creates).
$originallyActive.Activate()
